$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 297, shifting rows 297:399 down to 298:400
$ws.Rows("297:297").Insert()

# Populate the newly inserted row 297 with the new weekly data point.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R are constant across this consolidated
# dataset (same market/category), so we replicate them from the
# neighboring rows; D (date) and J (volume) carry the new unique values.
$ws.Cells.Item(297, 1).Value = 3
$ws.Cells.Item(297, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(297, 3).Value = "Coquimbo"
$ws.Cells.Item(297, 4).Value = 44985
$ws.Cells.Item(297, 5).Value = 5
$ws.Cells.Item(297, 6).Value = 100112039
$ws.Cells.Item(297, 7).Value = "Ciboulette"
$ws.Cells.Item(297, 8).Value = "Sin especificar"
$ws.Cells.Item(297, 9).Value = "Primera"
$ws.Cells.Item(297, 10).Value = 80
$ws.Cells.Item(297, 11).Value = 1500
$ws.Cells.Item(297, 12).Value = 1500
$ws.Cells.Item(297, 13).Value = 1500
$ws.Cells.Item(297, 14).Value = "`$/docena de atados"
$ws.Cells.Item(297, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(297, 16).Value = 500
$ws.Cells.Item(297, 17).Value = 3
$ws.Cells.Item(297, 18).Value = "Hortaliza"
